$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raw Results (In-Memory)")
$ws.Activate()

# Newly measured BSBM benchmark results for the early 0.4.2 build (rows 121-124, columns G-L)
$cols = @("G", "H", "I", "J", "K", "L")

$row121 = @(0.1643, 0.28999999999999998, 9.3620000000000001, 19225.830000000002, 0.18725, 0.18583)
$row122 = @(0.51680000000000004, 0.8448, 28.271999999999998, 6366.7, 0.56544000000000005, 0.56186999999999998)
$row123 = @(1.8153999999999999, 3.1198999999999999, 97.835999999999999, 1839.81, 1.95672, 1.9471499999999999)
$row124 = @(4.2595999999999998, 5.6853999999999996, 232.38200000000001, 774.59, 4.64764, 4.6351300000000002)

$rows = @{
    121 = $row121
    122 = $row122
    123 = $row123
    124 = $row124
}

foreach ($rowNum in $rows.Keys) {
    $values = $rows[$rowNum]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($rowNum, 7 + $i).Value = $values[$i]
    }
}

# Scroll the frozen pane down to reveal the newly-filled rows and move the
# active selection from G121 to G125, matching the updated view state
$window = $excel.ActiveWindow
$window.ScrollRow = 80
$ws.Range("G125").Select()
